# 965-MS-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-LateRepayment-Loanproduct.xlsx
# - Fix product-name label text (drop stray space: "Late Repayment" -> "LateRepayment")
#   on both the ProductLoanInput and ProductLoanOutput sheets.
# - Leave the workbook with the ProductLoanOutput sheet active/selected
#   (feature-file handling now drives verification from that sheet).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$newName = "965-MS-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-LateRepayment"

# Update the product-name value cell on the input sheet, then park the
# selection on B1 (matches the saved view state for this sheet).
$ws1.Range("B1").Value2 = $newName
$ws1.Range("B1").Select()

# Same value on the output sheet, then make it the active/selected sheet
# with its own selection parked on B1 as well.
$ws2.Range("B1").Value2 = $newName
$ws2.Activate()
$ws2.Range("B1").Select()
